$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "NUMBER" field-type label to "INTEGER" everywhere it is used
# (the field-type column, D, uses this label for several fields).
$usedRange = $ws.UsedRange
foreach ($cell in $usedRange.Cells) {
    $v = $cell.Value()
    if ($v -eq "NUMBER") {
        $cell.Value = "INTEGER"
    }
}

# Match the author's final cursor position noted in the edit.
$ws.Range("D16").Select()
